$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 694
    $ws.Range("F3").Value = 4022
    $ws.Range("F4").Value = 113
}
